$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster ECs -> Target cluster ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05881766666666666
$ws.Range("H2").Value = 0.176453
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.002172842242
$ws.Range("R2").Value = 0.019555580178
$ws.Range("S2").Value = 0.02099032928903418
$ws.Range("T2").Value = 0.02099032928903418

# Row 3 (Target cluster FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05881766666666666
$ws.Range("H3").Value = 0.176453
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("Q3").Value = 0.05547162763944444
$ws.Range("R3").Value = 0.499244648755
$ws.Range("S3").Value = 0.5358731102718634
$ws.Range("T3").Value = 0.5358731102718634

# Row 4 (Target cluster MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05881766666666666
$ws.Range("H4").Value = 0.176453
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("Q4").Value = 0.04587187862744444
$ws.Range("R4").Value = 0.412846907647
$ws.Range("S4").Value = 0.4431365604391025
$ws.Range("T4").Value = 0.4431365604391026
